$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 123, shifting existing rows 123-138 down to 124-139
$ws.Rows.Item(123).Insert()

# Fill in the constant columns (A-J) for the new row 123, matching the rest of the dataset
$ws.Cells.Item(123, 1).Value2 = 10
$ws.Cells.Item(123, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(123, 3).Value2 = "La Araucanía"
$ws.Cells.Item(123, 4).Value2 = 44516
$ws.Cells.Item(123, 5).Value2 = 9
$ws.Cells.Item(123, 6).Value2 = "Fruta"
$ws.Cells.Item(123, 7).Value2 = 100103
$ws.Cells.Item(123, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(123, 9).Value2 = 100103004
$ws.Cells.Item(123, 10).Value2 = "Durazno"

# Fill in the record-specific columns (K-T) for the new row
$ws.Cells.Item(123, 11).Value2 = "Early Majestic"
$ws.Cells.Item(123, 12).Value2 = "Primera"
$ws.Cells.Item(123, 13).Value2 = 65
$ws.Cells.Item(123, 14).Value2 = 18000
$ws.Cells.Item(123, 15).Value2 = 18000
$ws.Cells.Item(123, 16).Value2 = 18000
$ws.Cells.Item(123, 17).Value2 = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(123, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(123, 19).Value2 = 1800
$ws.Cells.Item(123, 20).Value2 = 10
